# Auto-generated edit script applying the Golem_Profits diff
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 250
$ws.Range("I9").Value = 250
$ws.Range("K9").Value = 250
$ws.Range("M9").Value = -81
$ws.Range("H21").Value = 800
$ws.Range("J21").Value = 800
$ws.Range("L21").Value = 800
$ws.Range("N21").Value = -1736
$ws.Range("H23").Value = 800
$ws.Range("J23").Value = 800
$ws.Range("L23").Value = 800
$ws.Range("N23").Value = -1268
$ws.Range("H38").Value = 557.6
$ws.Range("I38").Value = 431
$ws.Range("K38").Value = 1293
$ws.Range("M38").Value = -921
$ws.Range("H40").Value = 3479.8
$ws.Range("I40").Value = 1849.75
$ws.Range("K40").Value = 1849.75
$ws.Range("M40").Value = -1674.75
$ws.Range("H53").Value = 122.666664
$ws.Range("I53").Value = 85.28570999999999
$ws.Range("J53").Value = 175
$ws.Range("K53").Value = 85.28570999999999
$ws.Range("L53").Value = 175
$ws.Range("M53").Value = 551.71429
$ws.Range("N53").Value = -1449
$ws.Range("H100").Value = 539.1818
$ws.Range("J100").Value = 637.25
$ws.Range("L100").Value = 637.25
$ws.Range("N100").Value = -1719.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 10500
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2059.8572
$ws.Range("I105").Value = 2059.8572
$ws.Range("K105").Value = 2059.8572
$ws.Range("M105").Value = -312.8571999999999
$ws.Range("H107").Value = 1483.3793
$ws.Range("I107").Value = 1484.1818
$ws.Range("J107").Value = 1480.8572
$ws.Range("K107").Value = 1484.1818
$ws.Range("L107").Value = 1480.8572
$ws.Range("M107").Value = 435.8181999999999
$ws.Range("N107").Value = -5320.8572
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2373.0908
$ws.Range("I2").Value = 1242.8572
$ws.Range("J2").Value = 4351
$ws.Range("K2").Value = 1242.8572
$ws.Range("L2").Value = 4351
$ws.Range("M2").Value = -1129.8572
$ws.Range("N2").Value = -4577
$ws.Range("H50").Value = 43999
$ws.Range("J50").Value = 43999
$ws.Range("L50").Value = 43999
$ws.Range("N50").Value = -45249
$ws.Range("H58").Value = 1745
$ws.Range("I58").Value = 1745
$ws.Range("K58").Value = 1745
$ws.Range("M58").Value = -1542
$ws.Range("H86").Value = 8794.375
$ws.Range("I86").Value = 8980.714
$ws.Range("K86").Value = 8980.714
$ws.Range("M86").Value = -7857.714
$ws.Range("H89").Value = 8794.375
$ws.Range("I89").Value = 8980.714
$ws.Range("K89").Value = 44903.57
$ws.Range("M89").Value = -39287.57
$ws.Range("H93").Value = 23844.4
$ws.Range("I93").Value = 23844.4
$ws.Range("K93").Value = 23844.4
$ws.Range("M93").Value = -21972.4
$ws.Range("H107").Value = 394.4
$ws.Range("I107").Value = 394.4
$ws.Range("K107").Value = 394.4
$ws.Range("M107").Value = 1525.6
$ws.Range("H122").Value = 1096.6666
$ws.Range("I122").Value = 1151
$ws.Range("J122").Value = 988
$ws.Range("K122").Value = 3453
$ws.Range("L122").Value = 2964
$ws.Range("M122").Value = -1003
$ws.Range("N122").Value = -7864
$ws.Range("H136").Value = 1745
$ws.Range("I136").Value = 1745
$ws.Range("K136").Value = 5235
$ws.Range("M136").Value = -2685

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2511.1365
$ws.Range("I4").Value = 1874.6428
$ws.Range("J4").Value = 3625
$ws.Range("K4").Value = 5623.928400000001
$ws.Range("L4").Value = 10875
$ws.Range("M4").Value = -5511.928400000001
$ws.Range("N4").Value = -11099
$ws.Range("H14").Value = 5629.75
$ws.Range("I14").Value = 5629.75
$ws.Range("K14").Value = 16889.25
$ws.Range("M14").Value = -16716.25
$ws.Range("H86").Value = 2309.375
$ws.Range("I86").Value = 1166.6666
$ws.Range("J86").Value = 2995
$ws.Range("K86").Value = 3499.9998
$ws.Range("L86").Value = 8985
$ws.Range("M86").Value = -2313.9998
$ws.Range("N86").Value = -11357
$ws.Range("H89").Value = 2309.375
$ws.Range("I89").Value = 1166.6666
$ws.Range("J89").Value = 2995
$ws.Range("K89").Value = 10499.9994
$ws.Range("L89").Value = 26955
$ws.Range("M89").Value = -4571.999400000001
$ws.Range("N89").Value = -38811
$ws.Range("H107").Value = 703.3
$ws.Range("I107").Value = 366.66666
$ws.Range("J107").Value = 1208.25
$ws.Range("K107").Value = 1099.99998
$ws.Range("L107").Value = 3624.75
$ws.Range("M107").Value = 820.0000199999999
$ws.Range("N107").Value = -7464.75
$ws.Range("H114").Value = 2499.75
$ws.Range("I114").Value = 1666.3334
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 4999.0002
$ws.Range("L114").Value = 15000
$ws.Range("M114").Value = -1745.0002
$ws.Range("N114").Value = -21508
$ws.Range("H139").Value = 1378
$ws.Range("I139").Value = 1378
$ws.Range("K139").Value = 4134
$ws.Range("M139").Value = 1006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 45000
$ws.Range("J46").Value = 45000
$ws.Range("L46").Value = 45000
$ws.Range("N46").Value = -45312
$ws.Range("H62").Value = 65000
$ws.Range("J62").Value = 65000
$ws.Range("L62").Value = 65000
$ws.Range("N62").Value = -66372
$ws.Range("H65").Value = 65000
$ws.Range("J65").Value = 65000
$ws.Range("L65").Value = 195000
$ws.Range("N65").Value = -201864
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H97").Value = 449.6
$ws.Range("I97").Value = 437.25
$ws.Range("J97").Value = 499
$ws.Range("K97").Value = 437.25
$ws.Range("L97").Value = 499
$ws.Range("M97").Value = 58.75
$ws.Range("N97").Value = -1491
$ws.Range("H126").Value = 1000
$ws.Range("J126").Value = 1000
$ws.Range("L126").Value = 3000
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 2999.5
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 6000
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -3470
$ws.Range("N132").Value = -17057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -31322
$ws.Range("H132").Value = 3428.3333
$ws.Range("I132").Value = 2892.5
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 8677.5
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -6147.5
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H52").Value = 6000
$ws.Range("I52").Value = 6000
$ws.Range("K52").Value = 6000
$ws.Range("M52").Value = -5774
$ws.Range("H54").Value = 22757.223
$ws.Range("I54").Value = 21267.143
$ws.Range("J54").Value = 27972.5
$ws.Range("K54").Value = 21267.143
$ws.Range("L54").Value = 27972.5
$ws.Range("M54").Value = -20747.143
$ws.Range("N54").Value = -29012.5
$ws.Range("H107").Value = 1155.3636
$ws.Range("J107").Value = 1780.6
$ws.Range("L107").Value = 5341.799999999999
$ws.Range("N107").Value = -9181.799999999999
$ws.Range("H118").Value = 77000
$ws.Range("J118").Value = 77000
$ws.Range("L118").Value = 77000
$ws.Range("N118").Value = -80314
$ws.Range("H126").Value = 4744.4
$ws.Range("I126").Value = 4289
$ws.Range("K126").Value = 12867
$ws.Range("M126").Value = -10397
$ws.Range("H136").Value = 4981.4546
$ws.Range("I136").Value = 4981.4546
$ws.Range("K136").Value = 14944.3638
$ws.Range("M136").Value = -12394.3638
